# Updated WASM results and extended Emscripten build instructions in README
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated raw measurement data (columns H, I, J; rows 6-15)
$data = @(
    @(3.0752999999999999, 5.0811299999999999, 7.633),
    @(3.0070000000000001, 4.8851000000000004, 8.1645000000000003),
    @(3.3279999999999998, 4.7542,              7.5202),
    @(2.9969999999999999, 5.0053999999999998, 7.5201200000000004),
    @(3.0217000000000001, 4.9242999999999997, 7.8101000000000003),
    @(3.0537999999999998, 4.8785999999999996, 7.5704000000000002),
    @(2.9979,              4.9705000000000004, 7.5528000000000004),
    @(3.2791000000000001, 5.0410000000000004, 7.6109999999999998),
    @(3.2896000000000001, 5.0373999999999999, 7.6349999999999998),
    @(3.0366,              4.8520000000000003, 7.7354000000000003)
)

$row = 6
foreach ($vals in $data) {
    $ws.Cells.Item($row, 8).Value = $vals[0]
    $ws.Cells.Item($row, 9).Value = $vals[1]
    $ws.Cells.Item($row, 10).Value = $vals[2]
    $row++
}

# Update the active sheet view (scroll position, zoom, selection)
$window = $excel.ActiveWindow
$window.ScrollRow = 5
$window.ScrollColumn = 1
$window.Zoom = 205
$ws.Range("K21").Select()
